# Update cryptocurrency price and volume(1h) figures in the cryptos worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.501.57"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.362.16"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.73"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.40"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.361.04"
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.48"
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.935.76"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("E15").Value = "  -3.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.84"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.364.32"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.587.55"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.96"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.33"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "378.73"
$ws.Range("E22").Value = "  -2.92%  "
$ws.Range("E23").Value = "  -3.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.500.47"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.26"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("E28").Value = "  +10.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.53"
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.53"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E36").Value = "  -5.74%  "
$ws.Range("E37").Value = "  -2.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "165.78"
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("E40").Value = "  -3.93%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.769"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.50"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.39"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.23"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.83"
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.12"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.371.32"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  -2.50%  "